$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp header
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 19:50"

# Row 4
$ws.Range("B4").Value = 133146
$ws.Range("C4").Value = 9568
$ws.Range("E4").Value = 126405
$ws.Range("G4").Value = 143
$ws.Range("H4").Value = 2363

# Row 9
$ws.Range("A9").Value = "Francia"
$ws.Range("B9").Value = 40174
$ws.Range("C9").Value = 2599
$ws.Range("D9").Value = 7202
$ws.Range("E9").Value = 30366
$ws.Range("F9").Value = 4632
$ws.Range("G9").Value = 292
$ws.Range("H9").Value = 2606

# Row 10
$ws.Range("A10").Value = "Iran"
$ws.Range("B10").Value = 38309
$ws.Range("C10").Value = 2901
$ws.Range("D10").Value = 12391
$ws.Range("E10").Value = 23278
$ws.Range("F10").Value = 3206
$ws.Range("G10").Value = 123
$ws.Range("H10").Value = 2640

# Row 17
$ws.Range("B17").Value = 8711
$ws.Range("C17").Value = 440
$ws.Range("E17").Value = 8146

# Row 18
$ws.Range("A18").Value = "Canada"
$ws.Range("B18").Value = 6243
$ws.Range("C18").Value = 588
$ws.Range("D18").Value = 508
$ws.Range("E18").Value = 5672
$ws.Range("F18").Value = 120
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 63

# Row 19
$ws.Range("A19").Value = "Portugal"
$ws.Range("B19").Value = 5962
$ws.Range("C19").Value = 792
$ws.Range("D19").Value = 43
$ws.Range("E19").Value = 5800
$ws.Range("F19").Value = 89
$ws.Range("G19").Value = 19
$ws.Range("H19").Value = 119

# Row 21
$ws.Range("A21").Value = "Brasil"
$ws.Range("B21").Value = 4065
$ws.Range("C21").Value = 161
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 3941
$ws.Range("F21").Value = 296
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 118

# Row 22
$ws.Range("A22").Value = "Australia"
$ws.Range("B22").Value = 3980
$ws.Range("C22").Value = 345
$ws.Range("D22").Value = 226
$ws.Range("E22").Value = 3738
$ws.Range("F22").Value = 23
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 16

# Row 25
$ws.Range("E25").Value = 2716
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 16

# Row 44
$ws.Range("D44").Value = 95
$ws.Range("E44").Value = 902

# Row 58
$ws.Range("A58").Value = "Catar"
$ws.Range("B58").Value = 634
$ws.Range("C58").Value = 44
$ws.Range("D58").Value = 48
$ws.Range("E58").Value = 585
$ws.Range("F58").Value = 6
$ws.Range("H58").Value = 1

# Row 59
$ws.Range("A59").Value = "Colombia"
$ws.Range("B59").Value = 608
$ws.Range("D59").Value = 10
$ws.Range("E59").Value = 592
$ws.Range("F59").Value = 0
$ws.Range("H59").Value = 6

# Row 66
$ws.Range("B66").Value = 463
$ws.Range("C66").Value = 61
$ws.Range("E66").Value = 424

# Row 83
$ws.Range("A83").Value = "Jordania"
$ws.Range("B83").Value = 259
$ws.Range("C83").Value = 13
$ws.Range("D83").Value = 18
$ws.Range("E83").Value = 239
$ws.Range("F83").Value = 3
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 2

# Row 84
$ws.Range("A84").Value = "Kuwait"
$ws.Range("B84").Value = 255
$ws.Range("C84").Value = 20
$ws.Range("D84").Value = 67
$ws.Range("E84").Value = 188
$ws.Range("F84").Value = 12
$ws.Range("H84").Value = 0

# Row 85
$ws.Range("A85").Value = "Kazajistan"
$ws.Range("B85").Value = 251
$ws.Range("C85").Value = 23
$ws.Range("E85").Value = 232
$ws.Range("F85").Value = 0

# Row 144
$ws.Range("A144").Value = "Etiopia"
$ws.Range("B144").Value = 21
$ws.Range("C144").Value = 5
$ws.Range("D144").Value = 1
$ws.Range("E144").Value = 20

# Row 145
$ws.Range("A145").Value = "Congo"
$ws.Range("C145").Value = 15
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 19

# Row 175
$ws.Range("A175").Value = "Zimbabue"

# Row 176
$ws.Range("A176").Value = "Gabon"

# Row 177
$ws.Range("A177").Value = "Angola"
$ws.Range("B177").Value = 7
$ws.Range("C177").Value = 2
$ws.Range("E177").Value = 5
$ws.Range("G177").Value = 2
$ws.Range("H177").Value = 2

# Row 178
$ws.Range("A178").Value = "Santa Sede"

# Row 180
$ws.Range("A180").Value = "Eritrea"

# Row 181
$ws.Range("A181").Value = "Benin"
$ws.Range("E181").Value = 6
$ws.Range("H181").Value = 0

# Row 183
$ws.Range("A183").Value = "Cabo Verde"
$ws.Range("B183").Value = 6
$ws.Range("H183").Value = 1

# Row 184
$ws.Range("A184").Value = "Montserrat"
